# Improve the dictionary of nce2 articles
# Adds 5 new key/value rows (19-23) to the i18n dictionary sheet, reusing
# the same "key" style already used by the existing rows (copied via the
# row above so number formatting / borders / fill match exactly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows to append right after the existing data (row 18 is the last
# populated row: ToggleConciseMode / 切换简洁模式).
$newRows = @(
    @("TogglePanel", "显示/隐藏列表"),
    @("common",       "中文"),
    @("en",           "英语"),
    @("de",           "德语"),
    @("zh-CN",        "中文")
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Make sure the newly filled cells pick up the same style as the other
# key/value rows (border + fill + font used throughout A1:B18).
$srcRange = $ws.Range("A18:B18")
$dstRange = $ws.Range("A19:B23")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply the values since PasteSpecial only touched formatting, just to
# be safe in case any COM layer quirk cleared them.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
